$p = $ppt.ActivePresentation

# --- Slide 1: Subtitle "AKA / JSON vs GRPC vs SignalR Core" -> "AKA / REST vs GRPC vs SignalR Core"
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(4)
$tr1 = $subtitle.TextFrame.TextRange
$para1b = $tr1.Paragraphs(2)
$para1b.Text = "REST vs GRPC vs "
$para1b.InsertAfter("SignalR") | Out-Null
$para1b.InsertAfter(" Core") | Out-Null

# --- Slide 2: Content placeholder "JSON vs GRPC vs SignalR Core" -> "REST vs GRPC vs SignalR Core"
$slide2 = $p.Slides.Item(2)
$content2 = $slide2.Shapes.Item(5)
$tr2 = $content2.TextFrame.TextRange
$para2a = $tr2.Paragraphs(1)
$para2a.Text = "REST vs GRPC vs "
$para2a.InsertAfter("SignalR") | Out-Null
$para2a.InsertAfter(" Core") | Out-Null

# --- Slide 4: SmartArt diagram node "JSON with REST" -> "REST with JSON"
$slide4 = $p.Slides.Item(4)
$diagramShape = $slide4.Shapes.Item(4)
$smartArt = $diagramShape.SmartArt
$node = $smartArt.Nodes.Item(1)
$node.TextFrame.TextRange.Text = "REST with JSON"
